$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 388, shifting rows 388:477 down to 389:478
$ws.Rows.Item(388).Insert()

# Populate the new row 388 with data (copying constant columns from the row
# immediately below, which originally held the "old" row 388 data, and
# setting the new/changed values per the diff)
$ws.Range("A388").Value = 3
$ws.Range("B388").Value = "Femacal de La Calera"
$ws.Range("C388").Value = "Coquimbo"
$ws.Range("D388").Value = 44943
$ws.Range("E388").Value = 5
$ws.Range("F388").Value = 100114013
$ws.Range("G388").Value = "Zanahoria"
$ws.Range("H388").Value = "Sin especificar"
$ws.Range("I388").Value = "Primera"
$ws.Range("J388").Value = 480
$ws.Range("K388").Value = 11500
$ws.Range("L388").Value = 12000
$ws.Range("M388").Value = 11740
$ws.Range("N388").Value = "`$/saco 20 kilos"
$ws.Range("O388").Value = "Provincia de Quillota"
$ws.Range("P388").Value = 587
$ws.Range("Q388").Value = 20
$ws.Range("R388").Value = "Hortaliza"
